$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1213:H1213").Copy()
$ws.Range("A1233:H1233").PasteSpecial(-4122)
$ws.Range("A1233").Value = 43986
$ws.Range("B1233").Value = "Hospital Universitario Miguel Servet"
$ws.Range("C1233").Value = 7
$ws.Range("D1233").Value = 4
$ws.Range("E1233").Value = "Zaragoza"
$ws.Range("F1233").Value = "Zaragoza"
$ws.Range("G1233").Value = 50297
$ws.Range("H1233").Value = "Fuente Aragón Hoy"

$ws.Range("A1214:H1214").Copy()
$ws.Range("A1234:H1234").PasteSpecial(-4122)
$ws.Range("A1234").Value = 43986
$ws.Range("B1234").Value = "Hospital Clínico Universitario"
$ws.Range("C1234").Value = 27
$ws.Range("E1234").Value = "Zaragoza"
$ws.Range("F1234").Value = "Zaragoza"
$ws.Range("G1234").Value = 50297
$ws.Range("H1234").Value = "Fuente Aragón Hoy"

$ws.Range("A1215:H1215").Copy()
$ws.Range("A1235:H1235").PasteSpecial(-4122)
$ws.Range("A1235").Value = 43986
$ws.Range("B1235").Value = "Hospital Royo Villanova"
$ws.Range("C1235").Value = 4
$ws.Range("E1235").Value = "Zaragoza"
$ws.Range("F1235").Value = "Zaragoza"
$ws.Range("G1235").Value = 50297
$ws.Range("H1235").Value = "Fuente Aragón Hoy"

$ws.Range("A1216:H1216").Copy()
$ws.Range("A1236:H1236").PasteSpecial(-4122)
$ws.Range("A1236").Value = 43986
$ws.Range("B1236").Value = "Hospital Nuestra Señora de Gracia"
$ws.Range("C1236").Value = 1
$ws.Range("E1236").Value = "Zaragoza"
$ws.Range("F1236").Value = "Zaragoza"
$ws.Range("G1236").Value = 50297
$ws.Range("H1236").Value = "Fuente Aragón Hoy"

$ws.Range("A1217:H1217").Copy()
$ws.Range("A1237:H1237").PasteSpecial(-4122)
$ws.Range("A1237").Value = 43986
$ws.Range("B1237").Value = "Hospital General de la Defensa"
$ws.Range("C1237").Value = 1
$ws.Range("E1237").Value = "Zaragoza"
$ws.Range("F1237").Value = "Zaragoza"
$ws.Range("G1237").Value = 50297
$ws.Range("H1237").Value = "Fuente Aragón Hoy"

$ws.Range("A1218:H1218").Copy()
$ws.Range("A1238:H1238").PasteSpecial(-4122)
$ws.Range("A1238").Value = 43986
$ws.Range("B1238").Value = "Hospital Obispo Polanco"
$ws.Range("C1238").Value = 5
$ws.Range("E1238").Value = "Teruel"
$ws.Range("F1238").Value = "Teruel"
$ws.Range("G1238").Value = 44216
$ws.Range("H1238").Value = "Fuente Aragón Hoy"

$ws.Range("A1219:H1219").Copy()
$ws.Range("A1239:H1239").PasteSpecial(-4122)
$ws.Range("A1239").Value = 43986
$ws.Range("B1239").Value = "Hospital de Alcañiz"
$ws.Range("C1239").Value = 2
$ws.Range("E1239").Value = "Alcañiz"
$ws.Range("F1239").Value = "Teruel"
$ws.Range("G1239").Value = 44013
$ws.Range("H1239").Value = "Fuente Aragón Hoy"

$ws.Range("A1220:H1220").Copy()
$ws.Range("A1240:H1240").PasteSpecial(-4122)
$ws.Range("A1240").Value = 43986
$ws.Range("B1240").Value = "Hospital de Barbastro"
$ws.Range("C1240").Value = 11
$ws.Range("D1240").Value = 1
$ws.Range("E1240").Value = "Barbastro"
$ws.Range("F1240").Value = "Huesca"
$ws.Range("G1240").Value = 22048
$ws.Range("H1240").Value = "Fuente Aragón Hoy"

$ws.Range("A1221:H1221").Copy()
$ws.Range("A1241:H1241").PasteSpecial(-4122)
$ws.Range("A1241").Value = 43986
$ws.Range("B1241").Value = "Hospital San Jorge"
$ws.Range("C1241").Value = 7
$ws.Range("D1241").Value = 1
$ws.Range("E1241").Value = "Huesca"
$ws.Range("F1241").Value = "Huesca"
$ws.Range("G1241").Value = 22125
$ws.Range("H1241").Value = "Fuente Aragón Hoy"

$ws.Range("A1222:H1222").Copy()
$ws.Range("A1242:H1242").PasteSpecial(-4122)
$ws.Range("A1242").Value = 43986
$ws.Range("B1242").Value = "Hospital Sagrado Corazón"
$ws.Range("E1242").Value = "Huesca"
$ws.Range("F1242").Value = "Huesca"
$ws.Range("G1242").Value = 22125
$ws.Range("H1242").Value = "Fuente Aragón Hoy"

$ws.Range("A1223:H1223").Copy()
$ws.Range("A1243:H1243").PasteSpecial(-4122)
$ws.Range("A1243").Value = 43986
$ws.Range("B1243").Value = "Hospital Ernest Lluch"
$ws.Range("C1243").Value = 3
$ws.Range("E1243").Value = "Calatayud"
$ws.Range("F1243").Value = "Zaragoza"
$ws.Range("G1243").Value = 50067
$ws.Range("H1243").Value = "Fuente Aragón Hoy"

$ws.Range("A1224:H1224").Copy()
$ws.Range("A1244:H1244").PasteSpecial(-4122)
$ws.Range("A1244").Value = 43986
$ws.Range("B1244").Value = "Hospital San José"
$ws.Range("C1244").Value = 4
$ws.Range("E1244").Value = "Teruel"
$ws.Range("F1244").Value = "Teruel"
$ws.Range("G1244").Value = 44216
$ws.Range("H1244").Value = "Fuente Aragón Hoy"

$ws.Range("A1225:H1225").Copy()
$ws.Range("A1245:H1245").PasteSpecial(-4122)
$ws.Range("A1245").Value = 43986
$ws.Range("B1245").Value = "Hospital Ejea – Cinco Villas"
$ws.Range("E1245").Value = "Ejea de los Caballeros"
$ws.Range("F1245").Value = "Zaragoza"
$ws.Range("G1245").Value = 50095
$ws.Range("H1245").Value = "Fuente Aragón Hoy"

$ws.Range("A1226:H1226").Copy()
$ws.Range("A1246:H1246").PasteSpecial(-4122)
$ws.Range("A1246").Value = 43986
$ws.Range("B1246").Value = "MAZ"
$ws.Range("C1246").Value = 1
$ws.Range("E1246").Value = "Zaragoza"
$ws.Range("F1246").Value = "Zaragoza"
$ws.Range("G1246").Value = 50297
$ws.Range("H1246").Value = "Fuente Aragón Hoy"

$ws.Range("A1227:H1227").Copy()
$ws.Range("A1247:H1247").PasteSpecial(-4122)
$ws.Range("A1247").Value = 43986
$ws.Range("B1247").Value = "Hospital Viamed Montecanal"
$ws.Range("C1247").Value = 1
$ws.Range("E1247").Value = "Zaragoza"
$ws.Range("F1247").Value = "Zaragoza"
$ws.Range("G1247").Value = 50297
$ws.Range("H1247").Value = "Fuente Aragón Hoy"

$ws.Range("A1228:H1228").Copy()
$ws.Range("A1248:H1248").PasteSpecial(-4122)
$ws.Range("A1248").Value = 43986
$ws.Range("B1248").Value = "Clínica Montpellier"
$ws.Range("C1248").Value = 1
$ws.Range("E1248").Value = "Zaragoza"
$ws.Range("F1248").Value = "Zaragoza"
$ws.Range("G1248").Value = 50297
$ws.Range("H1248").Value = "Fuente Aragón Hoy"

$ws.Range("A1229:H1229").Copy()
$ws.Range("A1249:H1249").PasteSpecial(-4122)
$ws.Range("A1249").Value = 43986
$ws.Range("B1249").Value = "Hospital Quirón"
$ws.Range("C1249").Value = 1
$ws.Range("E1249").Value = "Zaragoza"
$ws.Range("F1249").Value = "Zaragoza"
$ws.Range("G1249").Value = 50297
$ws.Range("H1249").Value = "Fuente Aragón Hoy"

$ws.Range("A1230:H1230").Copy()
$ws.Range("A1250:H1250").PasteSpecial(-4122)
$ws.Range("A1250").Value = 43986
$ws.Range("B1250").Value = "Hospital San Juan de Dios de Zaragoza"
$ws.Range("E1250").Value = "Zaragoza"
$ws.Range("F1250").Value = "Zaragoza"
$ws.Range("G1250").Value = 50297
$ws.Range("H1250").Value = "Fuente Aragón Hoy"

$ws.Range("A1231:H1231").Copy()
$ws.Range("A1251:H1251").PasteSpecial(-4122)
$ws.Range("A1251").Value = 43986
$ws.Range("B1251").Value = "Clínica Viamed Santiago"
$ws.Range("E1251").Value = "Huesca"
$ws.Range("F1251").Value = "Huesca"
$ws.Range("G1251").Value = 22125
$ws.Range("H1251").Value = "Fuente Aragón Hoy"

$ws.Range("A1232:H1232").Copy()
$ws.Range("A1252:H1252").PasteSpecial(-4122)
$ws.Range("A1252").Value = 43986
$ws.Range("B1252").Value = "Clínica El Pilar"
$ws.Range("E1252").Value = "Zaragoza"
$ws.Range("F1252").Value = "Zaragoza"
$ws.Range("G1252").Value = 50297
$ws.Range("H1252").Value = "Fuente Aragón Hoy"
